$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Users" sheet: scores bumped up by 1 (Week 2 results applied)
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("C2").Value = 1   # Manuel
$wsUsers.Range("C3").Value = 2   # Daniel
$wsUsers.Range("C4").Value = 2   # Raff
$wsUsers.Range("C5").Value = 2   # Haunschi

# ---------------------------------------------------------------------------
# 2) "Matches" sheet: Week 2 matches (rows 18-30) get a Winner and are marked
#    Completed = TRUE.
# ---------------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches")

$wsMatches.Range("F18").Value = "Green Bay Packers"
$wsMatches.Range("G18").Value = $true

$wsMatches.Range("F19").Value = "Baltimore Ravens"
$wsMatches.Range("G19").Value = $true

$wsMatches.Range("F20").Value = "Cincinnati Bengals"
$wsMatches.Range("G20").Value = $true

$wsMatches.Range("F21").Value = "Dallas Cowboys"
$wsMatches.Range("G21").Value = $true

$wsMatches.Range("F22").Value = "Detroit Lions"
$wsMatches.Range("G22").Value = $true

$wsMatches.Range("F23").Value = "New England Patriots"
$wsMatches.Range("G23").Value = $true

$wsMatches.Range("F24").Value = "San Francisco 49ers"
$wsMatches.Range("G24").Value = $true

$wsMatches.Range("F25").Value = "Buffalo Bills"
$wsMatches.Range("G25").Value = $true

$wsMatches.Range("F26").Value = "Seattle Seahawks"
$wsMatches.Range("G26").Value = $true

$wsMatches.Range("F27").Value = "Los Angeles Rams"
$wsMatches.Range("G27").Value = $true

$wsMatches.Range("F28").Value = "Arizona Cardinals"
$wsMatches.Range("G28").Value = $true

$wsMatches.Range("F29").Value = "Indianapolis Colts"
$wsMatches.Range("G29").Value = $true

$wsMatches.Range("F30").Value = "Philadelphia Eagles"
$wsMatches.Range("G30").Value = $true

# ---------------------------------------------------------------------------
# 3) "Picks" sheet: add Week 2 picks (rows 6-9)
# ---------------------------------------------------------------------------
$wsPicks = $wb.Worksheets.Item("Picks")

$wsPicks.Range("A6").Value = 6
$wsPicks.Range("B6").Value = "Daniel"
$wsPicks.Range("C6").Value = 2
$wsPicks.Range("D6").Value = "Philadelphia Eagles @ Kansas City Chiefs"
$wsPicks.Range("E6").Value = "Philadelphia Eagles"
$wsPicks.Range("F6").Value = "Kansas City Chiefs"
$wsPicks.Range("G6").Value = "Philadelphia Eagles"
$wsPicks.Range("H6").Value = "YES"
$wsPicks.Range("I6").Value = 1

$wsPicks.Range("A7").Value = 8
$wsPicks.Range("B7").Value = "Haunschi"
$wsPicks.Range("C7").Value = 2
$wsPicks.Range("D7").Value = "Buffalo Bills @ New York Jets"
$wsPicks.Range("E7").Value = "Buffalo Bills"
$wsPicks.Range("F7").Value = "New York Jets"
$wsPicks.Range("G7").Value = "Buffalo Bills"
$wsPicks.Range("H7").Value = "YES"
$wsPicks.Range("I7").Value = 1

$wsPicks.Range("A8").Value = 5
$wsPicks.Range("B8").Value = "Manuel"
$wsPicks.Range("C8").Value = 2
$wsPicks.Range("D8").Value = "New York Giants @ Dallas Cowboys"
$wsPicks.Range("E8").Value = "Dallas Cowboys"
$wsPicks.Range("F8").Value = "New York Giants"
$wsPicks.Range("G8").Value = "Dallas Cowboys"
$wsPicks.Range("H8").Value = "YES"
$wsPicks.Range("I8").Value = 1

$wsPicks.Range("A9").Value = 7
$wsPicks.Range("B9").Value = "Raff"
$wsPicks.Range("C9").Value = 2
$wsPicks.Range("D9").Value = "New York Giants @ Dallas Cowboys"
$wsPicks.Range("E9").Value = "Dallas Cowboys"
$wsPicks.Range("F9").Value = "New York Giants"
$wsPicks.Range("G9").Value = "Dallas Cowboys"
$wsPicks.Range("H9").Value = "YES"
$wsPicks.Range("I9").Value = 1

# ---------------------------------------------------------------------------
# 4) "Team_Winner_Usage" sheet: rows 3-5 now reflect each user's *second*
#    used winning team (shifted from the prior Week-1-only contents), and
#    four new rows (6-9) record the Week-1 usage that got displaced.
# ---------------------------------------------------------------------------
$wsTWU = $wb.Worksheets.Item("Team_Winner_Usage")

$wsTWU.Range("A3").Value = "Daniel"
$wsTWU.Range("B3").Value = "Philadelphia Eagles"

$wsTWU.Range("A4").Value = "Haunschi"
$wsTWU.Range("B4").Value = "Buffalo Bills"

$wsTWU.Range("A5").Value = "Haunschi"
$wsTWU.Range("B5").Value = "Washington Commanders"

$wsTWU.Range("A6").Value = "Manuel"
$wsTWU.Range("B6").Value = "Atlanta Falcons"
$wsTWU.Range("C6").Value = 1
$wsTWU.Range("D6").Value = 2
$wsTWU.Range("E6").Value = "1/2x"
$wsTWU.Range("F6").Value = "YES"

$wsTWU.Range("A7").Value = "Manuel"
$wsTWU.Range("B7").Value = "Dallas Cowboys"
$wsTWU.Range("C7").Value = 1
$wsTWU.Range("D7").Value = 2
$wsTWU.Range("E7").Value = "1/2x"
$wsTWU.Range("F7").Value = "YES"

$wsTWU.Range("A8").Value = "Raff"
$wsTWU.Range("B8").Value = "Cincinnati Bengals"
$wsTWU.Range("C8").Value = 1
$wsTWU.Range("D8").Value = 2
$wsTWU.Range("E8").Value = "1/2x"
$wsTWU.Range("F8").Value = "YES"

$wsTWU.Range("A9").Value = "Raff"
$wsTWU.Range("B9").Value = "Dallas Cowboys"
$wsTWU.Range("C9").Value = 1
$wsTWU.Range("D9").Value = 2
$wsTWU.Range("E9").Value = "1/2x"
$wsTWU.Range("F9").Value = "YES"

# ---------------------------------------------------------------------------
# 5) "Summary" sheet: user score rows follow the same +1 bump as "Users"
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 1
$wsSummary.Range("D2").Value = "1 points"

$wsSummary.Range("C3").Value = 2
$wsSummary.Range("D3").Value = "2 points"

$wsSummary.Range("C4").Value = 2
$wsSummary.Range("D4").Value = "2 points"

$wsSummary.Range("C5").Value = 2
$wsSummary.Range("D5").Value = "2 points"
